$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "65.151.99"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.55%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.548.84"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "597.53"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "133.69"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.27%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.547.14"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("E12").Value = "  -1.23%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.155.60"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E14").Value = "  -3.31%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "26.89"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.86%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.549.99"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("E17").Value = "  -0.11%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "65.269.23"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "9.93"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.32%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.34"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  -1.12%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "390.11"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.88%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.576"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.85%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.692.35"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "74.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -1.10%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.77"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  +24.53%  "
$ws.Range("E30").Value = "  -0.12%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.53"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.76%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.07%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.551.27"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "24.04"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("E35").Value = "  -0.03%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.146"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.72%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "170.37"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("E39").Value = "  -0.93%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.02"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("E41").Value = "  -0.06%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.826"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.26%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "26.51"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "43.10"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("E45").Value = "  +4.51%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  -2.14%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.455.99"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.46%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.90"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0265"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
